{"js": "// Replace the arithmetic-problem text in the table cells according to the\n// mapping scraped from the diff. Each pair is a unique, exact-match\n// \"before\" string (e.g. \"723\u00d77=\") mapped to an exact \"after\" string\n// (e.g. \"942\u00d73=\"), so a case-sensitive, non-wildcard search + replace\n// is unambiguous for every cell (no \"before\" value repeats or collides\n// with another cell's \"after\" value).\nconst replacements = [\n  [\"723\u00d77=\", \"942\u00d73=\"],\n  [\"264\u00d75=\", \"911\u00d75=\"],\n  [\"719\u00d79=\", \"176\u00d72=\"],\n  [\"170\u00d72=\", \"784\u00d75=\"],\n  [\"905\u00d72=\", \"791\u00d72=\"],\n  [\"512\u00d78=\", \"850\u00d75=\"],\n  [\"447\u00d73=\", \"518\u00d75=\"],\n  [\"901\u00d74=\", \"112\u00d75=\"],\n  [\"367\u00d76=\", \"633\u00d77=\"],\n  [\"159\u00d77=\", \"630\u00d72=\"],\n  [\"712\u00d74=\", \"406\u00d73=\"],\n  [\"501\u00d79=\", \"184\u00d79=\"],\n  [\"597\u00d76=\", \"815\u00d73=\"],\n  [\"286\u00d72=\", \"306\u00d77=\"],\n  [\"197\u00d76=\", \"418\u00d77=\"],\n  [\"459\u00d75=\", \"527\u00d72=\"],\n  [\"869\u00d76=\", \"377\u00d74=\"],\n  [\"383\u00d77=\", \"268\u00d75=\"],\n  [\"252\u00d74=\", \"308\u00d76=\"],\n  [\"193\u00d78=\", \"186\u00d72=\"],\n  [\"482\u00d79=\", \"337\u00d72=\"],\n  [\"867\u00d77=\", \"598\u00d79=\"],\n  [\"647\u00d77=\", \"429\u00d77=\"],\n  [\"657\u00d77=\", \"407\u00d72=\"],\n  [\"313\u00d72=\", \"339\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const rng of results.items) {\n    rng.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the arithmetic-problem text in the table cells according to the\n# mapping scraped from the diff. Each \"before\" string is a unique, exact\n# token (e.g. \"723\u00d77=\") so Find/Replace \u2013 All over the whole document\n# story is unambiguous for every cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"723\u00d77=\", \"942\u00d73=\"),\n    @(\"264\u00d75=\", \"911\u00d75=\"),\n    @(\"719\u00d79=\", \"176\u00d72=\"),\n    @(\"170\u00d72=\", \"784\u00d75=\"),\n    @(\"905\u00d72=\", \"791\u00d72=\"),\n    @(\"512\u00d78=\", \"850\u00d75=\"),\n    @(\"447\u00d73=\", \"518\u00d75=\"),\n    @(\"901\u00d74=\", \"112\u00d75=\"),\n    @(\"367\u00d76=\", \"633\u00d77=\"),\n    @(\"159\u00d77=\", \"630\u00d72=\"),\n    @(\"712\u00d74=\", \"406\u00d73=\"),\n    @(\"501\u00d79=\", \"184\u00d79=\"),\n    @(\"597\u00d76=\", \"815\u00d73=\"),\n    @(\"286\u00d72=\", \"306\u00d77=\"),\n    @(\"197\u00d76=\", \"418\u00d77=\"),\n    @(\"459\u00d75=\", \"527\u00d72=\"),\n    @(\"869\u00d76=\", \"377\u00d74=\"),\n    @(\"383\u00d77=\", \"268\u00d75=\"),\n    @(\"252\u00d74=\", \"308\u00d76=\"),\n    @(\"193\u00d78=\", \"186\u00d72=\"),\n    @(\"482\u00d79=\", \"337\u00d72=\"),\n    @(\"867\u00d77=\", \"598\u00d79=\"),\n    @(\"647\u00d77=\", \"429\u00d77=\"),\n    @(\"657\u00d77=\", \"407\u00d72=\"),\n    @(\"313\u00d72=\", \"339\u00d79=\"),\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n# wdFindContinue = 1 (Wrap), wdReplaceAll = 2 (Replace)\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
